$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 with the new TPM-derived values (previously row 3's relationship,
# i.e. Resolving-Mac -> Resolving-Mac, recomputed), then delete the old row 3.
$ws.Range("D2").Value = "Resolving-Mac"
$ws.Range("G2").Value = 1.077228333333333
$ws.Range("H2").Value = 3.231685
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.1790523333333333
$ws.Range("N2").Value = 0.537157
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.1928802466161111
$ws.Range("R2").Value = 1.735922219545
$ws.Range("S2").Value = 1
$ws.Range("T2").Value = 1

# Remove the now-obsolete row 3 (the self-loop Resolving-Mac -> Resolving-Mac
# duplicate row that was merged into row 2 above).
$ws.Rows("3").Delete()
